$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$values = @(
    "холодно",
    "дороу",
    "привет",
    "дело",
    "спасибо",
    "сделать",
    "сказать",
    "сказать",
    "сделать",
    "дело",
    "сделать",
    "сделать",
    "сделать",
    "сказать",
    "дело",
    "распознавать",
    "распознавать",
    "распознавать",
    "распознавать",
    "распознавать",
    "привет",
    "дело",
    "дело",
    "спасибо",
    "дело",
    "дело",
    "дело",
    "погода",
    "спасибо",
    "хай",
    "распознавать",
    "сказать",
    "сказать",
    "распознавать",
    "распознавать",
    "распознавать",
    "сказать",
    "узнавать",
    "дело",
    "except",
    "погода норма",
    "привет",
    "распознавать",
    "распознавать",
    "стоять",
    "стоять",
    "садиться",
    "садиться"
)

$startRow = 18
for ($i = 0; $i -lt $values.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $values[$i]
}
